# chore: update Sheets via scheduled runner
#
# Recomputed Leve flip-profit figures (currentAveragePrice*, LevePrice*,
# LeveProfit* columns H:N) for a handful of rows across the job sheets
# after the upstream market-board pull. Values only -- no structural
# changes to the sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 32
$ws.Range("H32").Value = 810.2
$ws.Range("J32").Value = 862.75
$ws.Range("L32").Value = 862.75
$ws.Range("N32").Value = -1514.75
# row 55
$ws.Range("H55").Value = 304.45456
$ws.Range("J55").Value = 309.85715
$ws.Range("L55").Value = 309.85715
$ws.Range("N55").Value = -737.85715
# row 80
$ws.Range("H80").Value = 1169.4166
$ws.Range("I80").Value = 2173.8
$ws.Range("J80").Value = 452
$ws.Range("K80").Value = 6521.400000000001
$ws.Range("L80").Value = 1356
$ws.Range("M80").Value = -5523.400000000001
$ws.Range("N80").Value = -3352
# row 83
$ws.Range("H83").Value = 1169.4166
$ws.Range("I83").Value = 2173.8
$ws.Range("J83").Value = 452
$ws.Range("K83").Value = 19564.2
$ws.Range("L83").Value = 4068
$ws.Range("M83").Value = -14572.2
$ws.Range("N83").Value = -14052
# row 96
$ws.Range("H96").Value = 2230
$ws.Range("I96").Value = 750
$ws.Range("J96").Value = 2970
$ws.Range("K96").Value = 2250
$ws.Range("L96").Value = 8910
$ws.Range("M96").Value = -877
$ws.Range("N96").Value = -11656
# row 99
$ws.Range("H99").Value = 2297.75
$ws.Range("J99").Value = 2972.3333
$ws.Range("L99").Value = 8916.999899999999
$ws.Range("N99").Value = -11912.9999
# row 116
$ws.Range("H116").Value = 15098.625
$ws.Range("I116").Value = 34866.332
$ws.Range("K116").Value = 34866.332
$ws.Range("M116").Value = -31424.332
# row 132
$ws.Range("H132").Value = 1224.5
$ws.Range("I132").Value = 1120.0588
$ws.Range("K132").Value = 3360.1764
$ws.Range("M132").Value = -830.1764000000003
# row 137
$ws.Range("H137").Value = 1731
$ws.Range("I137").Value = 1357.4286
$ws.Range("K137").Value = 4072.2858
$ws.Range("M137").Value = -1522.2858
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 4127.021
$ws.Range("I32").Value = 2530.4324
$ws.Range("K32").Value = 2530.4324
$ws.Range("M32").Value = -2243.4324
# row 45
$ws.Range("H45").Value = 2082
$ws.Range("I45").Value = 1055.5714
$ws.Range("K45").Value = 1055.5714
$ws.Range("M45").Value = -678.5714
# row 74
$ws.Range("H74").Value = 3205.6
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 4676
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 4676
$ws.Range("M74").Value = -126
$ws.Range("N74").Value = -6424
# row 77
$ws.Range("H77").Value = 3205.6
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 4676
$ws.Range("K77").Value = 5000
$ws.Range("L77").Value = 23380
$ws.Range("M77").Value = -632
$ws.Range("N77").Value = -32116
# row 97
$ws.Range("H97").Value = 549.6667
$ws.Range("I97").Value = 549.6667
$ws.Range("K97").Value = 549.6667
$ws.Range("M97").Value = -53.66669999999999
# row 102
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622
# row 122
$ws.Range("H122").Value = 1661.5862
$ws.Range("I122").Value = 1637.9584
$ws.Range("K122").Value = 4913.8752
$ws.Range("M122").Value = -2463.8752
$ws = $wb.Worksheets.Item("BSM")
# row 108
$ws.Range("H108").Value = 34995.5
$ws.Range("J108").Value = 34995.5
$ws.Range("L108").Value = 34995.5
$ws.Range("N108").Value = -42675.5
$ws = $wb.Worksheets.Item("CUL")
# row 17
$ws.Range("H17").Value = 2979
$ws.Range("I17").Value = 2979
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 8937
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -8768
$ws.Range("N17").ClearContents()
# row 23
$ws.Range("H23").Value = 112.166664
$ws.Range("I23").Value = 45
$ws.Range("J23").Value = 145.75
$ws.Range("K23").Value = 135
$ws.Range("L23").Value = 437.25
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = -907.25
# row 34
$ws.Range("H34").Value = 1676
$ws.Range("J34").Value = 2493.3333
$ws.Range("L34").Value = 7479.999899999999
$ws.Range("N34").Value = -7647.999899999999
# row 38
$ws.Range("H38").Value = 283.68967
$ws.Range("I38").Value = 43
$ws.Range("J38").Value = 375.38095
$ws.Range("K38").Value = 129
$ws.Range("L38").Value = 1126.14285
$ws.Range("M38").Value = 218
$ws.Range("N38").Value = -1820.14285
# row 39
$ws.Range("H39").Value = 3883
$ws.Range("J39").Value = 3883
$ws.Range("L39").Value = 11649
$ws.Range("N39").Value = -12237
# row 55
$ws.Range("H55").Value = 18835.666
# row 80
$ws.Range("H80").Value = 5081.25
$ws.Range("I80").Value = 5108.3335
$ws.Range("K80").Value = 15325.0005
$ws.Range("M80").Value = -14389.0005
# row 83
$ws.Range("H83").Value = 5081.25
$ws.Range("I83").Value = 5108.3335
$ws.Range("K83").Value = 45975.0015
$ws.Range("M83").Value = -41295.0015
# row 131
$ws.Range("H131").Value = 11647203
$ws.Range("I131").Value = 83333910
$ws.Range("J131").Value = 22331.135
$ws.Range("K131").Value = 250001730
$ws.Range("L131").Value = 66993.405
$ws.Range("M131").Value = -249996690
$ws.Range("N131").Value = -77073.405
# row 132
$ws.Range("H132").Value = 1613.6666
$ws.Range("J132").Value = 1721.6
$ws.Range("L132").Value = 15494.4
$ws.Range("N132").Value = -20554.4
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 1463.4546
$ws.Range("I97").Value = 419.8
$ws.Range("J97").Value = 2333.1667
$ws.Range("K97").Value = 419.8
$ws.Range("L97").Value = 2333.1667
$ws.Range("M97").Value = 76.19999999999999
$ws.Range("N97").Value = -3325.1667
# row 113
$ws.Range("H113").Value = 1439.625
$ws.Range("I113").Value = 1241.75
$ws.Range("J113").Value = 1637.5
$ws.Range("K113").Value = 1241.75
$ws.Range("L113").Value = 1637.5
$ws.Range("M113").Value = 928.25
$ws.Range("N113").Value = -5977.5
# row 122
$ws.Range("H122").Value = 1998.409
$ws.Range("I122").Value = 1946.5834
$ws.Range("K122").Value = 5839.7502
$ws.Range("M122").Value = -3389.7502
# row 127
$ws.Range("H127").Value = 36639
$ws.Range("J127").Value = 36639
$ws.Range("L127").Value = 36639
$ws.Range("N127").Value = -46559
$ws = $wb.Worksheets.Item("LTW")
# row 55
$ws.Range("H55").Value = 307
$ws.Range("I55").Value = 95.888885
$ws.Range("K55").Value = 95.888885
$ws.Range("M55").Value = 77.111115
$ws = $wb.Worksheets.Item("WVR")
# row 48
$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 5000
$ws.Range("N48").Value = -6138
# row 100
$ws.Range("H100").Value = 250
$ws.Range("I100").Value = 250
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 500
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 41
$ws.Range("N100").ClearContents()
# row 123
$ws.Range("H123").Value = 46700
$ws.Range("J123").Value = 46700
$ws.Range("L123").Value = 46700
$ws.Range("N123").Value = -56500
